# The workbook tracks weekly Jengibre (ginger) price records for the
# "Terminal La Palmera de La Serena" market. A new weekly record is added
# at row 62 (just after the existing row 61), pushing all subsequent
# records down by one row. The sheet grows from 171 to 172 rows
# (170 data rows + header -> 171 data rows + header).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 62; this shifts rows 62..171 down to
# 63..172 and extends the used range / dimension to A1:R172 automatically.
$ws.Rows(62).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A62").Value = 8
$ws.Range("B62").Value = "Terminal La Palmera de La Serena"
$ws.Range("C62").Value = "Coquimbo"
$ws.Range("D62").Value = 45152
$ws.Range("E62").Value = 4
$ws.Range("F62").Value = 100114007
$ws.Range("G62").Value = "Jengibre"
$ws.Range("H62").Value = "Sin especificar"
$ws.Range("I62").Value = "Primera"
$ws.Range("J62").Value = 520
$ws.Range("K62").Value = 18000
$ws.Range("L62").Value = 19000
$ws.Range("M62").Value = 18500
$ws.Range("N62").Value = "`$/caja 13 kilos"
$ws.Range("O62").Value = "Perú"
$ws.Range("P62").Value = 1423
$ws.Range("Q62").Value = 13
$ws.Range("R62").Value = "Hortaliza"
